$d = $word.ActiveDocument

# Step 1: Replace paragraph containing "Building on the earlier findings..." with a
# clean single run (removing the spell-check proofErr wrapping around "XGBoost").
$paraA = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("Building on the earlier findings")) {
        $paraA = $cand
        break
    }
}
if ($paraA -eq $null) {
    throw "Could not locate 'Building on the earlier findings' paragraph"
}
$paraA.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Building on the earlier findings with the Random Forest model, the application of XGBoost to Bitcoin price prediction further underscores the significance of lagged features in time series forecasting. The comparative analysis reveals a stark contrast in the performance of the XGBoost model before and after the inclusion of lagged features. Initially, the XGBoost model—absent of these temporal indicators—reported a Mean Absolute Percentage Error (MAPE) of 25.99, a Mean Squared Prediction Error (MSPE) of 12.68, and an Akaike Information Criterion (AIC) of 34660.88, setting the stage for subsequent enhancement.</w:t></w:r></w:p>
'@)

# Step 2: Replace the paragraph containing "The integration of lagged features..." with a
# clean single run, and insert all of the new paragraphs (extra discussion paragraph,
# spacer paragraph with two line breaks, and the new "Model comparison:" section)
# immediately after it.
$paraC = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("The integration of lagged features into the")) {
        $paraC = $cand
        break
    }
}
if ($paraC -eq $null) {
    throw "Could not locate 'The integration of lagged features' paragraph"
}
$paraC.Range.InsertXML(@'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w14:ligatures w14:val="none"/></w:rPr><w:t>The integration of lagged features into the XGBoost model marks a pivotal improvement. The MAPE sees a dramatic reduction to 13.61, indicating a notable decrease in the percentage error of predictions. The MSPE follows suit, plummeting to 6.73, which signifies a greatly improved precision in forecasting the variance of Bitcoin prices. Additionally, the model's AIC drops to 27334.89, reflecting a refined model fit that better captures the complexities of the data.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w14:ligatures w14:val="none"/></w:rPr><w:t>This progression from the Random Forest to the XGBoost model, with the strategic incorporation of lagged features, demonstrates a consistent theme: temporal data points are invaluable for enhancing the accuracy and reliability of predictive models in financial time series. By effectively leveraging the information embedded in the preceding time steps, both models achieve a deeper level of analytical rigor, yielding forecasts that are not only more aligned with the actual market movements but also provide a stronger basis for decision-making in the volatile cryptocurrency domain.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="30"/><w:szCs w:val="30"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Model comparison:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="30"/><w:szCs w:val="30"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="30"/><w:szCs w:val="30"/><w14:ligatures w14:val="none"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w14:ligatures w14:val="none"/></w:rPr><w:t xml:space="preserve">The comparative analysis of Random Forest and XGBoost models for Bitcoin price prediction reveals a distinct enhancement in performance through the inclusion of lagged features. Initially, the Random Forest model presented a MAPE of 26.57 and an AIC of 34613.21. With lagged features, these metrics improved to a MAPE of 22.25 and an AIC of 27784.42. The XGBoost model showed even more substantial improvements post-lagged feature inclusion, with the MAPE dramatically reduced to 13.61 and the AIC to 27334.89. This comparison unequivocally </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w14:ligatures w14:val="none"/></w:rPr><w:lastRenderedPageBreak/><w:t>demonstrates the value of incorporating temporal data, with the XGBoost model achieving superior predictive accuracy and model fit over the Random Forest model. The integration of lagged features emerges as a pivotal factor in enhancing the models' ability to navigate the complexities of financial time series data.</w:t></w:r></w:p>
'@)

Write-Output "Edit applied successfully"
